# Weekly update: insert a new price-report row for "Zapallo italiano" at
# Vega Monumental Concepción, shifting all existing rows (95..151) down by
# one and appending the prior last row as the new final row (152).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 95; everything below (95..151)
# shifts down to (96..152).
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with this week's data.
$ws.Range("A95").Value = 11
$ws.Range("B95").Value = "Vega Monumental Concepción"
$ws.Range("C95").Value = "Bíobío"
$ws.Range("D95").Value = 44824
$ws.Range("E95").Value = 8
$ws.Range("F95").Value = 100112032
$ws.Range("G95").Value = "Zapallo italiano"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 100
$ws.Range("K95").Value = 16000
$ws.Range("L95").Value = 17000
$ws.Range("M95").Value = 16500
$ws.Range("N95").Value = '$/caja 50 unidades'
$ws.Range("O95").Value = "Región de Arica y Parinacota"
$ws.Range("P95").Value = 330
$ws.Range("Q95").Value = 50
$ws.Range("R95").Value = "Hortaliza"
